$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.346.93"
$ws.Range("E2").Value = "  -3.07%  "

$ws.Range("D3").Value = "1.746.65"
$ws.Range("E3").Value = "  -3.85%  "

$ws.Range("D4").Value = "1.009"
$ws.Range("E4").Value = "  +0.65%  "

$ws.Range("D5").Value = "321.67"
$ws.Range("E5").Value = "  -2.47%  "

$ws.Range("D6").Value = "1.007"
$ws.Range("E6").Value = "  +0.53%  "

$ws.Range("D7").Value = "0.4217"
$ws.Range("E7").Value = "  -4.57%  "

$ws.Range("D8").Value = "0.3605"
$ws.Range("E8").Value = "  -2.50%  "

$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").Value = "42.48"
$ws.Range("E9").Value = "  -4.83%  "

$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "0.07479"
$ws.Range("E10").Value = "  -2.75%  "

$ws.Range("D11").Value = "1.088"
$ws.Range("E11").Value = "  -3.31%  "

$ws.Range("D12").Value = "1.008"
$ws.Range("E12").Value = "  +0.67%  "

$ws.Range("D13").Value = "20.53"
$ws.Range("E13").Value = "  -6.80%  "

$ws.Range("D14").Value = "6.003"
$ws.Range("E14").Value = "  -4.07%  "

$ws.Range("D15").Value = "7.206"
$ws.Range("E15").Value = "  -4.50%  "

$ws.Range("D16").Value = "1.758.69"
$ws.Range("E16").Value = "  -3.32%  "

$ws.Range("D17").Value = "90.39"
$ws.Range("E17").Value = "  -2.50%  "

$ws.Range("D18").Value = "0.00001067"
$ws.Range("E18").Value = "  -1.56%  "

$ws.Range("D19").Value = "0.06354"
$ws.Range("E19").Value = "  -3.84%  "

$ws.Range("D20").Value = "1.005"
$ws.Range("E20").Value = "  +0.40%  "

$ws.Range("D21").Value = "16.92"
$ws.Range("E21").Value = "  -3.48%  "

$ws.Range("D22").Value = "5.863"
$ws.Range("E22").Value = "  -5.39%  "

$ws.Range("D23").Value = "27.446.82"
$ws.Range("E23").Value = "  -2.86%  "

$ws.Range("D24").Value = "11.08"
$ws.Range("E24").Value = "  -5.10%  "

$ws.Range("D25").Value = "2.084"
$ws.Range("E25").Value = "  +2.49%  "

$ws.Range("D26").Value = "161.14"
$ws.Range("E26").Value = "  +3.31%  "

$ws.Range("D27").Value = "20.13"
$ws.Range("E27").Value = "  -2.89%  "

$ws.Range("D28").Value = "1.989.87"
$ws.Range("E28").Value = "  -1.59%  "

$ws.Range("D29").Value = "2.109"
$ws.Range("E29").Value = "  -8.87%  "

$ws.Range("D30").Value = "123.89"
$ws.Range("E30").Value = "  -3.33%  "

$ws.Range("D31").Value = "1.098"
$ws.Range("E31").Value = "  -8.68%  "

$ws.Range("D32").Value = "3.659"
$ws.Range("E32").Value = "  -0.34%  "

$ws.Range("D33").Value = "5.509"
$ws.Range("E33").Value = "  -6.03%  "

$ws.Range("D34").Value = "0.08853"
$ws.Range("E34").Value = "  -3.88%  "

$ws.Range("D35").Value = "12.18"
$ws.Range("E35").Value = "  -6.53%  "

$ws.Range("D36").Value = "0.02272"
$ws.Range("E36").Value = "  -3.45%  "

$ws.Range("D37").Value = "0.2088"
$ws.Range("E37").Value = "  -3.77%  "

$ws.Range("D38").Value = "0.05993"
$ws.Range("E38").Value = "  -3.56%  "

$ws.Range("D39").Value = "0.6293"
$ws.Range("E39").Value = "  -4.07%  "

$ws.Range("D40").Value = "4.905"
$ws.Range("E40").Value = "  -4.68%  "

$ws.Range("D41").Value = "1.179"
$ws.Range("E41").Value = "  -1.51%  "

$ws.Range("D42").Value = "1.006"
$ws.Range("E42").Value = "  +0.51%  "

$ws.Range("D43").Value = "1.396"
$ws.Range("E43").Value = "  +0.60%  "

$ws.Range("D44").Value = "7.836"
$ws.Range("E44").Value = "  -3.86%  "

$ws.Range("D45").Value = "13.31"
$ws.Range("E45").Value = "  -3.36%  "

$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D46").Value = "3.688"
$ws.Range("E46").Value = "  -2.04%  "

$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "0.5835"
$ws.Range("E47").Value = "  -3.93%  "

$ws.Range("D48").Value = "122.77"
$ws.Range("E48").Value = "  -3.44%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.970"
$ws.Range("E49").Value = "  -3.14%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.150"
$ws.Range("E50").Value = "  -0.24%  "

$ws.Range("D51").Value = "0.06809"
$ws.Range("E51").Value = "  -2.33%  "
